$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "(15 points): As a user, I want to be able to randomly re-select ..."
#    Strike through the whole paragraph and highlight every run yellow.
# ---------------------------------------------------------------------------
$paraA = $d.Paragraphs.Item(14)
$paraA.Range.Font.StrikeThrough = 1
$paraA.Range.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 2) "(10 points): As a user, I want to be able to confirm that my day trip
#    is "complete" ..." - same treatment (strike + yellow highlight).
# ---------------------------------------------------------------------------
$paraB = $d.Paragraphs.Item(15)
$paraB.Range.Font.StrikeThrough = 1
$paraB.Range.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 3) "(5 points): As a developer, I want all of my functions to have a
#    Single Responsibility. Remember, each function should do just one
#    thing!" - the first three runs ("(5 points): ", "As a developer, I
#    want ", "all of my functions to have a Single Responsibility") collapse
#    into one run; every run in the paragraph (incl. the paragraph mark)
#    gets struck through. No highlight here.
# ---------------------------------------------------------------------------

# Re-saving identical text over the found range forces Word to rebuild the
# run(s) it spans into a single run - this is how the three leading runs get
# merged into one.
$mergeText = "(5 points): As a developer, I want all of my functions to have a Single Responsibility"
$mergeRng = $d.Content
$mergeRng.Find.Execute($mergeText, $false, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2) | Out-Null

# Walk forward run-by-run, striking each piece through. Setting
# Font.StrikeThrough on a precise sub-range splits a run at that boundary
# without merging it back into its (identically formatted) neighbours, so
# this reproduces the exact run layout the diff expects.
$cur = $d.Content
$cur.Find.Execute($mergeText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1
$pos = $cur.End

$cur = $d.Range($pos, $d.Content.End)
$cur.Find.Execute(".", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1
$pos = $cur.End

$cur = $d.Range($pos, $d.Content.End)
$cur.Find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1
$pos = $cur.End

$cur = $d.Range($pos, $d.Content.End)
$cur.Find.Execute("R", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1
$pos = $cur.End

$cur = $d.Range($pos, $d.Content.End)
$cur.Find.Execute("emember, each function should do just one thing", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1
$pos = $cur.End

$cur = $d.Range($pos, $d.Content.End)
$cur.Find.Execute("!", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cur.Font.StrikeThrough = 1

# Also strike the paragraph mark itself (w:pPr/w:rPr) so the paragraph-level
# formatting matches, without disturbing the run splits made above. Go back
# through the document's own Paragraphs collection (not a sub-range) so the
# full paragraph extent - not just the last "!" - is affected.
$paraC = $d.Paragraphs.Item(17)
$paraC.Range.Font.StrikeThrough = 1
